$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 1.53
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 5.75
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("W9").Value = 8.5
$ws.Range("AF9").Value = 41
$ws.Range("AH9").Value = 19
$ws.Range("AK9").Value = 51
$ws.Range("AW9").Value = 7
$ws.Range("AY9").Value = 29
